$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(338, 1).Value = "TestA"
$ws.Cells.Item(338, 1).ClearFormats()
